$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.696.13"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "3.173.37"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'535.34"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").Value = "'142.60"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.172.00"
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("D9").Value = "'0.449"
$ws.Range("E9").Value = "  +2.64%  "
$ws.Range("D10").Value = "'7.25"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").Value = "'0.110"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").Value = "'0.400"
$ws.Range("E12").Value = "  +4.42%  "
$ws.Range("D13").Value = "3.726.06"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("E14").Value = "  +2.89%  "
$ws.Range("D15").Value = "'25.94"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("D17").Value = "58.796.95"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "3.183.90"
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("D19").Value = "'6.19"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").Value = "'12.94"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "'8.09"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "'357.29"
$ws.Range("E22").Value = "  +5.37%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'0.515"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "'68.41"
$ws.Range("E25").Value = "  +3.23%  "
$ws.Range("D26").Value = "'0.170"
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("D27").Value = "0.0₃0961"
$ws.Range("E27").Value = "  +5.86%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "'7.53"
$ws.Range("E29").Value = "  +4.10%  "
$ws.Range("D30").Value = "'6.56"
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D33").Value = "'21.38"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").Value = "'1.22"
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("E35").Value = "  +7.06%  "
$ws.Range("D36").Value = "'157.83"
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("D37").Value = "'6.23"
$ws.Range("E37").Value = "  +3.20%  "
$ws.Range("D38").Value = "'26.53"
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("D40").Value = "'1.66"
$ws.Range("E40").Value = "  +13.62%  "
$ws.Range("D41").Value = "'0.0678"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D42").Value = "'0.708"
$ws.Range("E42").Value = "  +4.09%  "
$ws.Range("D43").Value = "'4.05"
$ws.Range("E43").Value = "  +4.10%  "
$ws.Range("D44").Value = "3.217.45"
$ws.Range("E44").Value = "  +2.20%  "
$ws.Range("D45").Value = "'0.0275"
$ws.Range("E45").Value = "  +7.29%  "
$ws.Range("D46").Value = "'36.94"
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("D47").Value = "2.346.72"
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").Value = "'1.02"
$ws.Range("E49").Value = "  +6.12%  "
$ws.Range("D50").Value = "'20.70"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "'6.09"
$ws.Range("E51").Value = "  +1.39%  "
